$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.175
$ws.Range("C2").Value = 0.5928571428571429
$ws.Range("J2").Value = 0.003571428571428571
$ws.Range("P2").Value = 0.1035714285714286
$ws.Range("S2").Value = 0.125
$ws.Range("C3").Value = 0.005952380952380952
$ws.Range("J3").Value = 0.01785714285714286
$ws.Range("P3").Value = 0.8214285714285714
$ws.Range("S3").Value = 0.1547619047619048
$ws.Range("J4").Value = 0.0196078431372549
$ws.Range("P4").Value = 0.6862745098039216
$ws.Range("S4").Value = 0.2941176470588235
$ws.Range("P5").Value = 0.75
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.07623318385650224
$ws.Range("D6").Value = 0.0179372197309417
$ws.Range("E6").Value = 0.004484304932735426
$ws.Range("F6").Value = 0.09417040358744394
$ws.Range("J6").Value = 0.273542600896861
$ws.Range("O6").Value = 0.0179372197309417
$ws.Range("Q6").Value = 0.1838565022421525
$ws.Range("R6").Value = 0.01345291479820628
$ws.Range("S6").Value = 0.3183856502242152
$ws.Range("B7").Value = 0.0972972972972973
$ws.Range("D7").Value = 0.02702702702702703
$ws.Range("E7").Value = 0.005405405405405406
$ws.Range("F7").Value = 0.03243243243243243
$ws.Range("J7").Value = 0.1135135135135135
$ws.Range("O7").Value = 0.01621621621621622
$ws.Range("Q7").Value = 0.2162162162162162
$ws.Range("R7").Value = 0.03243243243243243
$ws.Range("S7").Value = 0.4594594594594595
$ws.Range("B8").Value = 0.09554140127388536
$ws.Range("D8").Value = 0.02547770700636943
$ws.Range("E8").Value = 0.002123142250530786
$ws.Range("F8").Value = 0.05095541401273886
$ws.Range("J8").Value = 0.1061571125265393
$ws.Range("O8").Value = 0.02760084925690021
$ws.Range("Q8").Value = 0.1549893842887473
$ws.Range("R8").Value = 0.08492569002123142
$ws.Range("S8").Value = 0.4522292993630573
$ws.Range("B9").Value = 0.1003861003861004
$ws.Range("D9").Value = 0.0193050193050193
$ws.Range("F9").Value = 0.06177606177606178
$ws.Range("J9").Value = 0.111969111969112
$ws.Range("O9").Value = 0.0193050193050193
$ws.Range("Q9").Value = 0.1814671814671815
$ws.Range("R9").Value = 0.08880308880308881
$ws.Range("S9").Value = 0.416988416988417
$ws.Range("B10").Value = 0.1020746887966805
$ws.Range("D10").Value = 0.02157676348547718
$ws.Range("E10").Value = 0.002489626556016597
$ws.Range("F10").Value = 0.06473029045643154
$ws.Range("J10").Value = 0.1286307053941909
$ws.Range("O10").Value = 0.01659751037344398
$ws.Range("Q10").Value = 0.1867219917012448
$ws.Range("R10").Value = 0.07053941908713693
$ws.Range("S10").Value = 0.4066390041493776
$ws.Range("G11").Value = 0.147887323943662
$ws.Range("J11").Value = 0.102112676056338
$ws.Range("K11").Value = 0.2147887323943662
$ws.Range("L11").Value = 0.5176056338028169
$ws.Range("S11").Value = 0.0176056338028169
$ws.Range("G12").Value = 0.7870967741935484
$ws.Range("J12").Value = 0.1548387096774194
$ws.Range("L12").Value = 0.006451612903225806
$ws.Range("S12").Value = 0.05161290322580645
$ws.Range("G13").Value = 0.7368421052631579
$ws.Range("J13").Value = 0.131578947368421
$ws.Range("S13").Value = 0.131578947368421
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.00423728813559322
$ws.Range("H15").Value = 0.1567796610169492
$ws.Range("I15").Value = 0.09745762711864407
$ws.Range("J15").Value = 0.3432203389830508
$ws.Range("K15").Value = 0.03389830508474576
$ws.Range("M15").Value = 0.008474576271186441
$ws.Range("O15").Value = 0.0635593220338983
$ws.Range("S15").Value = 0.2923728813559322
$ws.Range("F16").Value = 0.01005025125628141
$ws.Range("H16").Value = 0.1206030150753769
$ws.Range("I16").Value = 0.1206030150753769
$ws.Range("J16").Value = 0.457286432160804
$ws.Range("K16").Value = 0.09547738693467336
$ws.Range("M16").Value = 0.01507537688442211
$ws.Range("O16").Value = 0.04522613065326633
$ws.Range("S16").Value = 0.135678391959799
$ws.Range("F17").Value = 0.02142857142857143
$ws.Range("H17").Value = 0.1380952380952381
$ws.Range("I17").Value = 0.1285714285714286
$ws.Range("J17").Value = 0.3785714285714286
$ws.Range("K17").Value = 0.09047619047619047
$ws.Range("M17").Value = 0.01904761904761905
$ws.Range("O17").Value = 0.08333333333333333
$ws.Range("S17").Value = 0.1404761904761905
$ws.Range("F18").Value = 0.0189873417721519
$ws.Range("H18").Value = 0.1265822784810127
$ws.Range("I18").Value = 0.120253164556962
$ws.Range("J18").Value = 0.3354430379746836
$ws.Range("K18").Value = 0.1518987341772152
$ws.Range("M18").Value = 0.0189873417721519
$ws.Range("O18").Value = 0.0759493670886076
$ws.Range("S18").Value = 0.1518987341772152
$ws.Range("F19").Value = 0.02631578947368421
$ws.Range("H19").Value = 0.2368421052631579
$ws.Range("I19").Value = 0.1009957325746799
$ws.Range("J19").Value = 0.3271692745376956
$ws.Range("K19").Value = 0.09530583214793741
$ws.Range("M19").Value = 0.01635846372688478
$ws.Range("N19").Value = 0.0007112375533428165
$ws.Range("O19").Value = 0.06401137980085349
$ws.Range("S19").Value = 0.1322901849217639
